$d = $word.ActiveDocument

# Disable smart-quote autoformatting so replacement text keeps straight quotes
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# Update the "Last compiled on" date: December 20, 2021 -> January 16, 2022
$d.Content.Find.Execute("December", $true, $false, $false, $false, $false, $true, 1, $false, "January", 2)
$d.Content.Find.Execute("20,", $true, $false, $false, $false, $false, $true, 1, $false, "16,", 2)
$d.Content.Find.Execute("2021", $true, $false, $false, $false, $false, $true, 1, $false, "2022", 2)

# Heading case fix: Load Libraries -> Load libraries
$d.Content.Find.Execute("Load Libraries", $true, $false, $false, $false, $false, $true, 1, $false, "Load libraries", 2)

# Dataset filename update (set Range.Text directly so the straight quotes
# around the file name are preserved, rather than letting Find/Replace's
# smart-quote autoformat turn them into curly quotes)
$rQuoted = $d.Content
$rQuoted.Find.Execute("""nhanes_homework.rda""")
$rQuoted.Text = """nhanes_homework_dataset.rda"""
